$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Mask the mobile numbers in column C (rows 2-5) with a placeholder string
$ws.Range("C2").Value = "98********"
$ws.Range("C3").Value = "98********"
$ws.Range("C4").Value = "98********"
$ws.Range("C5").Value = "98********"

# Update the active selection from B5 to C5
$ws.Range("C5").Select()
